$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the previously-entered exam grade values (students haven't taken Exam 2 yet)
$ws.Range("D12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C21").ClearContents()

# Update selection to reflect where the review of Exam 2 grading started
# (student scrolled down and selected D12:D14, ending on D14)
$ws.Range("D12:D14").Select()
